$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records arrived (date 2023-10-04 / serial 45203) for the
# "Chirimoya" subset at Macroferia Regional de Talca. They belong at the top
# of the historical block (row 52), so push the existing rows 52:160 down by
# two rows before writing the new data in.
$ws.Rows("52:53").Insert()

# New row 52: Primera
$ws.Range("A52").Value = 5
$ws.Range("B52").Value = "Macroferia Regional de Talca"
$ws.Range("C52").Value = "Maule"
$ws.Range("D52").Value = 45203
$ws.Range("E52").Value = 7
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100107
$ws.Range("H52").Value = "Otros"
$ws.Range("I52").Value = 100107002
$ws.Range("J52").Value = "Chirimoya"
$ws.Range("K52").Value = "Cultivar IV Región"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 100
$ws.Range("N52").Value = 24000
$ws.Range("O52").Value = 24000
$ws.Range("P52").Value = 24000
$ws.Range("Q52").Value = "$/bandeja 10 kilos"
$ws.Range("R52").Value = "Provincia de Limarí"
$ws.Range("S52").Value = 2400
$ws.Range("T52").Value = 10

# New row 53: Segunda
$ws.Range("A53").Value = 5
$ws.Range("B53").Value = "Macroferia Regional de Talca"
$ws.Range("C53").Value = "Maule"
$ws.Range("D53").Value = 45203
$ws.Range("E53").Value = 7
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100107
$ws.Range("H53").Value = "Otros"
$ws.Range("I53").Value = 100107002
$ws.Range("J53").Value = "Chirimoya"
$ws.Range("K53").Value = "Cultivar IV Región"
$ws.Range("L53").Value = "Segunda"
$ws.Range("M53").Value = 180
$ws.Range("N53").Value = 20000
$ws.Range("O53").Value = 20000
$ws.Range("P53").Value = 20000
$ws.Range("Q53").Value = "$/bandeja 10 kilos"
$ws.Range("R53").Value = "Provincia de Limarí"
$ws.Range("S53").Value = 2000
$ws.Range("T53").Value = 10
